$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.384.69"
$ws.Range("E2").Value = "  -5.75%  "
$ws.Range("D3").Value = "1.834.49"
$ws.Range("E3").Value = "  -4.64%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4217"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -8.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07207"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8977"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.76%  "
$ws.Range("D12").Value = "1.914.30"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.554"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.317"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06811"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "77.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008940"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.87%  "
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -8.45%  "
$ws.Range("D21").Value = "27.417.84"
$ws.Range("E21").Value = "  -5.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.927"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.29%  "
$ws.Range("D24").Value = "2.049.20"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.022"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.74%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("E27").Value = "  -4.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.222"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "110.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.681"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08851"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7706"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.491"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.857"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  -14.21%  "
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.084"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01923"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.786"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5025"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.62%  "
$ws.Range("E43").Value = "  -7.24%  "
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.159"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4697"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "104.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.20%  "
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("E49").Value = "  -9.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.630"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.830"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -16.11%  "
